# Minor change to .Rmd for HTML
#
# Extends the closing paragraph of the practical with a pointer to the
# "R for Data Science" online book's chapter on R Markdown, ending with a
# hyperlink to https://r4ds.had.co.nz/r-markdown.html

$d = $word.ActiveDocument

$oldText = "Now that you have the repository locally, start to play around with editing settings, headings, lines in either this practical, or the git practical, to get a good grasp of the capabilities."

$linkUrl = "https://r4ds.had.co.nz/r-markdown.html"

$newText = "Now that you have the repository locally, start to play around with editing settings, headings, lines in either this practical, or the git practical, to get a good grasp of the capabilities. There is an excellent one-page summary of the main things you can do with RMarkdown in the " + [char]0x201C + "R for Data Science" + [char]0x201D + " book, available online at " + $linkUrl

# 1) Grow the existing sentence with the new lead-in text, ending with the
#    plain-text URL (so it can be located again below).
$found = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

if (-not $found) {
    throw "Could not find the target paragraph text to replace."
}

# 2) Locate the plain-text URL we just inserted and turn it into a real
#    hyperlink, so the document ends up with a <w:hyperlink> run exactly
#    like the other links in this document.
$target = $d.Content.Duplicate
$target.Find.Execute($linkUrl, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$link = $d.Hyperlinks.Add($target, $linkUrl)
